$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Financial affidavit" entry belongs alphabetically between "Fee waiver" (row 15)
# and "Housing discrimination complaint - IDHR" (old row 16). Insert a new row at 16;
# this shifts rows 16:39 down to 17:40, carrying their existing cell values along.
$ws.Rows("16:16").Insert()

# Populate the newly inserted row.
$ws.Range("A16").Value = "Financial affidavit"
$ws.Range("B16").Value = "https://www.illinoislegalaid.org/legal-information/financial-affidavit"

# NOTE: Row/column insert in this runtime moves cell values but does not re-anchor
# existing Hyperlink objects to their new cells, so every hyperlink whose row was
# pushed down needs to be re-added at its new address (its URL/target is unchanged).
$ws.Hyperlinks.Add($ws.Range("B35"), "https://www.illinoislegalaid.org/legal-information/security-deposit-demand-letter")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://www.illinoislegalaid.org/legal-information/housing-discrimination-complaint-idhr")
$ws.Hyperlinks.Add($ws.Range("B38"), "https://www.illinoislegalaid.org/legal-information/stop-wage-assignment-letter")
$ws.Hyperlinks.Add($ws.Range("B31"), "https://www.illinoislegalaid.org/legal-information/request-time-work-due-domestic-abuse-letter")
$ws.Hyperlinks.Add($ws.Range("B32"), "https://www.illinoislegalaid.org/legal-information/respond-lawsuit")
$ws.Hyperlinks.Add($ws.Range("B40"), "https://www.illinoislegalaid.org/legal-information/voluntary-acknowledgment-parentage-vap")
$ws.Hyperlinks.Add($ws.Range("B18"), "https://www.illinoislegalaid.org/legal-information/interpreter-request")
$ws.Hyperlinks.Add($ws.Range("B20"), "https://www.illinoislegalaid.org/legal-information/motion")
$ws.Hyperlinks.Add($ws.Range("B39"), "https://www.illinoislegalaid.org/legal-information/transfer-death-instrument-or-todi")
$ws.Hyperlinks.Add($ws.Range("B27"), "https://www.illinoislegalaid.org/legal-information/power-attorney-agent-resign-letter")
$ws.Hyperlinks.Add($ws.Range("B28"), "https://www.illinoislegalaid.org/legal-information/power-attorney-revocation")
$ws.Hyperlinks.Add($ws.Range("B26"), "https://www.illinoislegalaid.org/legal-information/power-attorney-property")
$ws.Hyperlinks.Add($ws.Range("B25"), "https://www.illinoislegalaid.org/legal-information/power-attorney-health-care")
$ws.Hyperlinks.Add($ws.Range("B23"), "https://www.illinoislegalaid.org/legal-information/order-protection")
$ws.Hyperlinks.Add($ws.Range("B21"), "https://www.illinoislegalaid.org/legal-information/name-change-adult")
$ws.Hyperlinks.Add($ws.Range("B36"), "https://www.illinoislegalaid.org/legal-information/short-term-guardian-appointment")
$ws.Hyperlinks.Add($ws.Range("B30"), "https://www.illinoislegalaid.org/legal-information/remove-eviction-public-record")
$ws.Hyperlinks.Add($ws.Range("B33"), "https://www.illinoislegalaid.org/legal-information/respond-eviction")
$ws.Hyperlinks.Add($ws.Range("B37"), "https://www.illinoislegalaid.org/legal-information/small-claims-complaint")
$ws.Hyperlinks.Add($ws.Range("B34"), "https://www.illinoislegalaid.org/legal-information/security-deposit-complaint")

# New hyperlink for the inserted row itself.
$ws.Hyperlinks.Add($ws.Range("B16"), "https://www.illinoislegalaid.org/legal-information/financial-affidavit")

# Restore the selection/active cell as recorded after the edit.
$ws.Range("C16").Select()
